$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($ws, $row, $values)
    $cols = @("B","C","D","E","F","G","H","I","J")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}

# Sheet 1: Edi06_2_bg_detlim
$ws1 = $wb.Worksheets.Item("Edi06_2_bg_detlim")

Set-RowValues $ws1 2 @(0.023, 0.022, 0.023, 0.02, 0.021, 0.022, 0.001, 0.02, 0.023)
Set-RowValues $ws1 5 @(0.023, 0.022, 0.023, 0.02, 0.021, 0.022, 0.001, 0.02, 0.023)
Set-RowValues $ws1 8 @(0.036, 0.034, 0.036, 0.031, 0.033, 0.034, 0.002, 0.031, 0.036)

# Sheet 2: Edi06_3_bg_apf_detlim
$ws2 = $wb.Worksheets.Item("Edi06_3_bg_apf_detlim")

Set-RowValues $ws2 2 @(0.028, 0.026, 0.028, 0.024, 0.025, 0.026, 0.002, 0.024, 0.028)
Set-RowValues $ws2 5 @(0.028, 0.026, 0.028, 0.024, 0.025, 0.026, 0.002, 0.024, 0.028)
Set-RowValues $ws2 8 @(0.044, 0.041, 0.043, 0.038, 0.04, 0.041, 0.003, 0.038, 0.044)
